$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge the banner ranges that are going to gain a new column J. ---
$ws.Range("C1:I2").UnMerge()
$ws.Range("A3:I3").UnMerge()

# --- 2. Stamp column J with the same look as column I for every banner
#        row BEFORE re-merging, so the merge has nothing left to stamp. ---
$ws.Range("I1:I3").Copy()
$ws.Range("J1:J3").PasteSpecial(-4122)
$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("C1:J2").Merge()
$ws.Range("A3:J3").Merge()

# --- 3. Unmerge the old footer ranges and create the new (empty / unstyled)
#        footer merges BEFORE applying the saved format, so Merge() has
#        nothing pre-existing to stamp across the range. ---
$ws.Range("A17:D17").UnMerge()
$ws.Range("E17:I17").UnMerge()

$ws.Range("F13:J13").Merge()
$ws.Range("E17").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("A13:E13").Merge()

# --- 5. Update the filter/banner line (row 4). ---
$ws.Range("A4").Value = "Depot : opiant"
$ws.Range("C4").Value = "From : 01-10-2018"
$ws.Range("E4").Value = "To : 24-04-2019"

# --- 6. Update the table header row (row 7). ---
$ws.Range("H7").Value = "Audited"
$ws.Range("I7").Value = "Audited By"
$ws.Range("J7").Value = "Cash Remitted By"

# --- 7. Clear out the old data rows 8-10 entirely, then rewrite rows 8-11. ---
$ws.Range("A8:K10").Clear()

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 55237734
$ws.Range("D8").Value = "O2-2-1-Afternoon"
$ws.Range("F8").Value = "subhash(123123)"
$ws.Range("G8").Value = "RJ27BE4554"
$ws.Range("H8").Value = "Un-audited"

$ws.Range("A9").Value = 2
$ws.Range("B9").Value = 55237733
$ws.Range("D9").Value = "O2-2-1-Afternoon"
$ws.Range("F9").Value = "Apurv(242424)"
$ws.Range("G9").Value = "up13ba2296"
$ws.Range("H9").Value = "Un-audited"

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 55237734
$ws.Range("C10").Value = "01-03-2019 15:34:34"
$ws.Range("D10").Value = "O2-2-1-Afternoon"
$ws.Range("F10").Value = "Apurv(242424)"
$ws.Range("G10").Value = "RJ27BE4554"
$ws.Range("H10").Value = "Audited"
$ws.Range("I10").Value = "Satya"
$ws.Range("J10").Value = "Satya"

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = 55237734
$ws.Range("D11").Value = "O2-2-1-Afternoon"
$ws.Range("F11").Value = "Apurv(242424)"
$ws.Range("G11").Value = "RJ27BE4554"
$ws.Range("H11").Value = "Un-audited"

# --- 8. Write the footer row text, then wipe the old row 17 entirely. ---
$ws.Range("A13").Value = "Print taken by : Satya"
$ws.Range("F13").Value = "Print taken at : 24-04-2019 14:52:42"

$ws.Range("A17:K17").Clear()

# --- 9. Fix selection / active cell. ---
$ws.Range("F13").Select()
